$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '39.579.99'
Set-TextValue "E2" '  +1.82%  '

# Row 3
Set-TextValue "D3" '2.159.91'
Set-TextValue "E3" '  +1.71%  '

# Row 4
Set-TextValue "E4" '  -0.05%  '

# Row 5
Set-TextValue "D5" '226.90'
Set-TextValue "E5" '  -0.61%  '

# Row 6
Set-TextValue "D6" '0.622'
Set-TextValue "E6" '  +0.80%  '

# Row 7
Set-TextValue "D7" '62.80'
Set-TextValue "E7" '  +0.81%  '

# Row 8
Set-TextValue "E8" '  +0.01%  '

# Row 9
Set-TextValue "D9" '0.390'
Set-TextValue "E9" '  -0.02%  '

# Row 10
Set-TextValue "D10" '0.0843'
Set-TextValue "E10" '  -0.22%  '

# Row 11
Set-TextValue "D11" '0.104'
Set-TextValue "E11" '  +0.56%  '

# Row 12
Set-TextValue "D12" '15.85'
Set-TextValue "E12" '  -0.49%  '

# Row 13
Set-TextValue "D13" '2.479.98'
Set-TextValue "E13" '  +1.78%  '

# Row 14
Set-TextValue "D14" '21.72'
Set-TextValue "E14" '  -1.72%  '

# Row 15
Set-TextValue "D15" '0.803'
Set-TextValue "E15" '  -0.66%  '

# Row 16
Set-TextValue "D16" '5.47'
Set-TextValue "E16" '  -0.67%  '

# Row 17
Set-TextValue "D17" '2.151.67'
Set-TextValue "E17" '  +2.23%  '

# Row 18
Set-TextValue "D18" '39.573.38'
Set-TextValue "E18" '  +1.82%  '

# Row 19
Set-TextValue "D19" '71.58'
Set-TextValue "E19" '  -0.30%  '

# Row 20
Set-TextValue "D20" '6.04'
Set-TextValue "E20" '  -0.75%  '

# Row 21
Set-TextValue "D21" '0.0₃0854'
Set-TextValue "E21" '  +0.92%  '

# Row 22
Set-TextValue "D22" '227.76'
Set-TextValue "E22" '  +0.01%  '

# Row 23
Set-TextValue "E23" '  +0.04%  '

# Row 24
Set-TextValue "E24" '  +1.17%  '

# Row 25
Set-TextValue "D25" '2.28'
Set-TextValue "E25" '  -3.63%  '

# Row 26
Set-TextValue "D26" '170.44'
Set-TextValue "E26" '  +0.13%  '

# Row 27
Set-TextValue "D27" '9.43'
Set-TextValue "E27" '  -1.13%  '

# Row 28
Set-TextValue "E28" '  +1.10%  '

# Row 29
Set-TextValue "D29" '1.43'
Set-TextValue "E29" '  +0.63%  '

# Row 30
Set-TextValue "D30" '19.62'
Set-TextValue "E30" '  +0.86%  '

# Row 31
Set-TextValue "D31" '2.68'
Set-TextValue "E31" '  +4.53%  '

# Row 32
Set-TextValue "D32" '0.122'
Set-TextValue "E32" '  +0.60%  '

# Row 33
Set-TextValue "D33" '4.56'
Set-TextValue "E33" '  -0.74%  '

# Row 34
Set-TextValue "D34" '4.71'
Set-TextValue "E34" '  -1.69%  '

# Row 35
Set-TextValue "D35" '6.97'
Set-TextValue "E35" '  -3.24%  '

# Row 36
Set-TextValue "D36" '0.0617'
Set-TextValue "E36" '  +0.19%  '

# Row 37
Set-TextValue "D37" '3.82'
Set-TextValue "E37" '  +7.81%  '

# Row 38
Set-TextValue "D38" '2.39'
Set-TextValue "E38" '  -0.18%  '

# Row 39
Set-TextValue "D39" '5.09'
Set-TextValue "E39" '  +22.82%  '

# Row 40
Set-TextValue "E40" '  -0.11%  '

# Row 41
Set-TextValue "D41" '102.51'
Set-TextValue "E41" '  -0.01%  '

# Row 42
Set-TextValue "E42" '  -0.91%  '

# Row 43
Set-TextValue "D43" '17.63'
Set-TextValue "E43" '  -2.90%  '

# Row 44
Set-TextValue "D44" '1.514.20'
Set-TextValue "E44" '  -0.98%  '

# Row 45
Set-TextValue "E45" '  +0.09%  '

# Row 46
Set-TextValue "D46" '7.86'
Set-TextValue "E46" '  +0.98%  '

# Row 47
Set-TextValue "B47" 'HuobiToken'
Set-TextValue "C47" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D47" '2.80'
Set-TextValue "E47" '  +0.13%  '

# Row 48
Set-TextValue "B48" 'Cronos'
Set-TextValue "C48" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D48" '0.0917'
Set-TextValue "E48" '  -0.15%  '

# Row 49
Set-TextValue "E49" '  -0.35%  '

# Row 50
Set-TextValue "D50" '0.000191'
Set-TextValue "E50" '  +27.43%  '

# Row 51
Set-TextValue "D51" '2.99'
Set-TextValue "E51" '  +0.64%  '
